$d = $word.ActiveDocument

# --- 1. First paragraph: append trailing spaces to the existing text, then
#        add three red-colored runs forming "(This is a change – Version for main branch)"
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$dash = [char]0x2013
$p1 = $d.Paragraphs(1).Range
$insertPoint = $p1.End - 1

$run1Text = "(This is a change " + $dash + " Ve"
$run2Text = "rsion for main branch"
$run3Text = ")"

$rng1 = $d.Range($insertPoint, $insertPoint)
$rng1.InsertAfter($run1Text)
$rng1.Font.Color = 255

$insertPoint = $insertPoint + $run1Text.Length
$rng2 = $d.Range($insertPoint, $insertPoint)
$rng2.InsertAfter($run2Text)
$rng2.Font.Color = 255

$insertPoint = $insertPoint + $run2Text.Length
$rng3 = $d.Range($insertPoint, $insertPoint)
$rng3.InsertAfter($run3Text)
$rng3.Font.Color = 255

# --- 2. Append a new, otherwise-empty paragraph after the last paragraph in
#        the document, shaded with fill color F9F9F9.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$shadedParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($shadedParaXml)
